$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("courses")

# The promotion's "department" cell used to hold the school name; this
# batch update repurposes it to the generic "Packages" label.
$ws.Range("C2").Value = "Packages"

# Match the author's saved selection state.
$ws.Range("J23").Select()
